$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @{ 2=1.02; 3=1.028325096075421; 4=1.032745402505437; 5=1.028305576243729; 6=1.040452075360148; 9=1.036409555177593; 10=1.033478041777262; 11=1.035549747146059; 12=1.031122778701696; 13=1.043234384260475; 14=1.034945699138184 }
    3 = @{ 2=1.02; 3=1.029167659468975; 4=1.033385438919104; 5=1.029018398859814; 6=1.042653199654112; 9=1.03668972602798; 10=1.033961714541885; 11=1.035999071181374; 12=1.031643765603641; 13=1.045242280792426; 14=1.035430058773601 }
    4 = @{ 2=1.02; 3=1.02971298304797; 4=1.033799532065413; 5=1.029480134258528; 6=1.044072019043702; 9=1.036869357086506; 10=1.034274165066734; 11=1.036289074287637; 12=1.031980699621381; 13=1.046535740374956; 14=1.035742953014034 }
    5 = @{ 2=1.02; 3=1.029942267880508; 4=1.033973604101237; 5=1.029674365052266; 6=1.044667216360309; 9=1.036944478418706; 10=1.034405395454395; 11=1.036410814972321; 12=1.032122303702204; 13=1.047078155598337; 14=1.035874369763897 }
    6 = @{ 2=1.02; 3=1.02998076761802; 4=1.034002830815852; 5=1.029706984118683; 6=1.044767078791116; 9=1.036957068465367; 10=1.034427422359454; 11=1.036431245424067; 12=1.032146077128391; 13=1.047169150880603; 14=1.035896427949688 }
    7 = @{ 2=1.02; 3=1.029716046642168; 4=1.033801858074512; 5=1.029482729119853; 6=1.044079977068253; 9=1.036870362414023; 10=1.034275919059871; 11=1.036290701686134; 12=1.031982591911025; 13=1.046542993442886; 14=1.035744709498042 }
    8 = @{ 2=1.02; 3=1.028609817663396; 4=1.032961716844407; 5=1.028546375961717; 6=1.041197104322817; 9=1.036504584434108; 10=1.033641608667309; 11=1.035701751716021; 12=1.031298885778723; 13=1.043914177699312; 14=1.035109498311966 }
    9 = @{ 2=1.02; 3=1.026661485727244; 4=1.031480870439162; 5=1.026900184035261; 6=1.036073900951829; 9=1.035847271934023; 10=1.032519893351259; 11=1.034658261016774; 12=1.030092735158811; 13=1.039236237673968; 14=1.033986190031505 }
    10 = @{ 2=1.02; 3=1.025363253949684; 4=1.030493355455034; 5=1.025805285786077; 6=1.032627357489214; 9=1.035400387316523; 10=1.031769388296655; 11=1.033958744687122; 12=1.029287708733598; 13=1.036085049554788; 14=1.033234619173639 }
    11 = @{ 2=1.02; 3=1.024801258150795; 4=1.030065679987918; 5=1.025331793822606; 6=1.031127141164438; 9=1.035204801659914; 10=1.031443766774793; 11=1.033654923376144; 12=1.028938902043167; 13=1.034712403721044; 14=1.032908535231867 }
    12 = @{ 2=1.02; 3=1.024592529622534; 4=1.029906810615171; 5=1.025156009051365; 6=1.030568676655248; 9=1.035131837722027; 10=1.031322718542131; 11=1.033541930641791; 12=1.028809305663347; 13=1.034201279124083; 14=1.032787315096831 }
    13 = @{ 2=1.02; 3=1.024637301623156; 4=1.029940889169494; 5=1.025193711351246; 6=1.030688524917372; 9=1.035147503008086; 10=1.031348688238577; 11=1.033566174314638; 12=1.028837106071652; 13=1.034310974827525; 14=1.032813321673224 }
    14 = @{ 2=1.02; 3=1.024784004144421; 4=1.030052548036369; 5=1.025317261536836; 6=1.031081003378581; 9=1.035198776869732; 10=1.031433762890655; 11=1.033645586221803; 12=1.028928190262033; 13=1.034670179978588; 14=1.032898517141065 }
    15 = @{ 2=1.02; 3=1.02487439527794; 4=1.030121343191469; 5=1.025393396976166; 6=1.031322659806569; 9=1.03523032662542; 10=1.031486167201329; 11=1.033694495945581; 12=1.028984305718677; 13=1.034891329680533; 14=1.032950995871872 }
    16 = @{ 2=1.02; 3=1.02540055482895; 4=1.030521737254701; 5=1.02583672267392; 6=1.032726753367977; 9=1.035413323667008; 10=1.031790985027106; 11=1.033978888724356; 12=1.029310853125042; 13=1.036175972478853; 14=1.033256246573924 }
    17 = @{ 2=1.02; 3=1.025730639953246; 4=1.030772873532855; 5=1.026114971585824; 6=1.03360537745915; 9=1.035527554178866; 10=1.031982015351281; 11=1.034157032274959; 12=1.029515627296296; 13=1.03697958437876; 14=1.033447548183072 }
    18 = @{ 2=1.02; 3=1.025923187245369; 4=1.030919349872444; 5=1.026277327951093; 6=1.034117110253443; 9=1.035593982188918; 10=1.032093377636946; 11=1.034260851095435; 12=1.029635046922184; 13=1.037447532741739; 14=1.033559068615956 }
    19 = @{ 2=1.02; 3=1.02598884331892; 4=1.030969293308365; 5=1.026332697190154; 6=1.034291471243098; 9=1.035616598437677; 10=1.032131338718712; 11=1.034296235508035; 12=1.029675762257398; 13=1.037606959018144; 14=1.033597083606812 }
    20 = @{ 2=1.02; 3=1.025695223474863; 4=1.030745929734862; 5=1.026085112060754; 6=1.033511187659292; 9=1.035515319101088; 10=1.031961526071118; 11=1.034137928400522; 12=1.029493659210385; 13=1.036893445873585; 14=1.033427029805781 }
    21 = @{ 2=1.02; 3=1.024740803311533; 4=1.030019667608437; 5=1.025280876567243; 6=1.030965462179608; 9=1.03518368669921; 10=1.031408713233289; 11=1.033622205265204; 12=1.028901369180127; 13=1.034564438173019; 14=1.032873431910314 }
    22 = @{ 2=1.02; 3=1.024140847014316; 4=1.029562970083469; 5=1.024775749650022; 6=1.029357800167402; 9=1.034973354382848; 10=1.031060571131971; 11=1.033297139388029; 12=1.02852877581402; 13=1.033092773638531; 14=1.032524795407279 }
    23 = @{ 2=1.02; 3=1.02445888343088; 4=1.029805080643435; 5=1.025043476909007; 6=1.030210735095156; 9=1.035085028884491; 10=1.031245181762842; 11=1.033469540051638; 12=1.028726313382422; 13=1.03387363762092; 14=1.032709668206421 }
    24 = @{ 2=1.02; 3=1.025711226613459; 4=1.030758104499988; 5=1.026098604114536; 6=1.033553750303107; 9=1.035520848226363; 10=1.03197078448816; 11=1.034146560895218; 12=1.029503585710092; 13=1.036932370580068; 14=1.033436301370837 }
    25 = @{ 2=1.02; 3=1.027165060032425; 4=1.031863754059688; 5=1.027325313732162; 6=1.037403698819391; 9=1.036018725202053; 10=1.032810357111142; 11=1.034928705447357; 12=1.030404716916377; 13=1.04045119664797; 14=1.034277066283243 }
}

foreach ($row in $data.Keys) {
    $cols = $data[$row]
    foreach ($col in $cols.Keys) {
        $ws.Cells.Item($row, $col).Value = $cols[$col]
    }
}
